# Apply commit edit: refresh job-listing rows with a new LinkedIn search
# pull, add four new tracking columns (O-R), and trim the sheet from 26
# down to 24 rows (2 stale listings removed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new columns O-R + updated "Checked companies" count ---
$ws.Range("O1").Value = "Search Keyword"
$ws.Range("P1").Value = "Search Country"
$ws.Range("Q1").Value = "Job Date"
$ws.Range("R1").Value = "transformed publish date from description"
$ws.Range("T1").Value = 23

# Give the four new header cells the same bold/border/centered look as the
# existing header row (copy formatting only, values already set above).
$ws.Range("A1").Copy()
$ws.Range("O1:R1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column R holds plain-text ISO dates ("2025-04-17", ...); format as Text
# first so Excel does not auto-convert them into date serials.
$ws.Range("R2:R24").NumberFormat = "@"

# --- Row data (rows 2-24): one 1x18 array per row, written in a single
#     Range(...).Value assignment (columns A:R) ---
# Row 2: Sayva Solutions - Financial Planning and Analysis Manager
$row2 = New-Object "object[,]" 1,18
$row2[0,0] = "Sayva Solutions"
$row2[0,1] = "Financial Planning and Analysis Manager"
$row2[0,2] = $false
$row2[0,3] = $false
$row2[0,4] = $false
$row2[0,5] = $true
$row2[0,6] = $false
$row2[0,7] = $true
$row2[0,8] = $false
$row2[0,9] = $false
$row2[0,10] = "https://www.linkedin.com/jobs/view/4283551489/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=gpaeOmaMcT6ITyN7ucr0uA%3D%3D&trackingId=gRQQKomMl8RmekK8J5SHDg%3D%3D&trk=flagship3_search_srp_jobs"
$row2[0,11] = 92.35
$row2[0,12] = "financial, about, erp, present, lead, accounting, experience, years, business, key"
$row2[0,13] = "forecasting, remote"
$row2[0,14] = "remote job, planning manager"
$row2[0,15] = "United States"
$row2[0,16] = "5 days ago"
$row2[0,17] = "2025-04-17"
$ws.Range("A2:R2").Value = $row2

# Row 3: Flexton Inc. - Senior Project Manager
$row3 = New-Object "object[,]" 1,18
$row3[0,0] = "Flexton Inc."
$row3[0,1] = "Senior Project Manager"
$row3[0,2] = $false
$row3[0,3] = $false
$row3[0,4] = $false
$row3[0,5] = $true
$row3[0,6] = $false
$row3[0,7] = $true
$row3[0,8] = $false
$row3[0,9] = $false
$row3[0,10] = "https://www.linkedin.com/jobs/view/4286692400/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=gpaeOmaMcT6ITyN7ucr0uA%3D%3D&trackingId=522JsYcu3wxo3pCat60UdA%3D%3D&trk=flagship3_search_srp_jobs"
$row3[0,11] = 115.49
$row3[0,12] = "labor, project, forecasting, lead, manager, integrations, soft, new, est, stakeholder"
$row3[0,13] = "forecasting, remote"
$row3[0,14] = "remote job, planning manager"
$row3[0,15] = "United States"
$row3[0,16] = "2 days ago"
$row3[0,17] = "2025-04-20"
$ws.Range("A3:R3").Value = $row3

# Row 4: UGG - Sr. Director, Demand Planning - UGG
$row4 = New-Object "object[,]" 1,18
$row4[0,0] = "UGG"
$row4[0,1] = "Sr. Director, Demand Planning - UGG"
$row4[0,2] = $false
$row4[0,3] = $false
$row4[0,4] = $false
$row4[0,5] = $true
$row4[0,6] = $false
$row4[0,7] = $true
$row4[0,8] = $false
$row4[0,9] = $false
$row4[0,10] = "https://www.linkedin.com/jobs/view/4281738674/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=9v7sPhR2I6v10s%2B%2B%2FJrI5w%3D%3D&trackingId=MEy9S1l0llvZANwkh8Z2lA%3D%3D&trk=flagship3_search_srp_jobs"
$row4[0,11] = 176.46
$row4[0,12] = "plan, demand, planning, get, experience, channel, time, work, employee, market"
$row4[0,13] = "demand planning, forecasting, remote"
$row4[0,14] = "remote job, planning manager"
$row4[0,15] = "United States"
$row4[0,16] = "4 days ago"
$row4[0,17] = "2025-04-18"
$ws.Range("A4:R4").Value = $row4

# Row 5: The Sage Group - Regional Operations and Project Manager
$row5 = New-Object "object[,]" 1,18
$row5[0,0] = "The Sage Group"
$row5[0,1] = "Regional Operations and Project Manager"
$row5[0,2] = $false
$row5[0,3] = $false
$row5[0,4] = $false
$row5[0,5] = $true
$row5[0,6] = $false
$row5[0,7] = $true
$row5[0,8] = $false
$row5[0,9] = $false
$row5[0,10] = "https://www.linkedin.com/jobs/view/4272437959/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=9v7sPhR2I6v10s%2B%2B%2FJrI5w%3D%3D&trackingId=kPWMKiYQwAHja1vEBc9pBg%3D%3D&trk=flagship3_search_srp_jobs"
$row5[0,11] = 179.09
$row5[0,12] = "manage, project, management, community, site, lead, property, experience, operational, housing"
$row5[0,13] = "forecasting, remote"
$row5[0,14] = "remote job, planning manager"
$row5[0,15] = "United States"
$row5[0,16] = "3 weeks ago"
$row5[0,17] = "2025-04-01"
$ws.Range("A5:R5").Value = $row5

# Row 6: UGG - Sr. Director, Demand Planning - UGG
$row6 = New-Object "object[,]" 1,18
$row6[0,0] = "UGG"
$row6[0,1] = "Sr. Director, Demand Planning - UGG"
$row6[0,2] = $false
$row6[0,3] = $false
$row6[0,4] = $false
$row6[0,5] = $true
$row6[0,6] = $false
$row6[0,7] = $true
$row6[0,8] = $false
$row6[0,9] = $false
$row6[0,10] = "https://www.linkedin.com/jobs/view/4281744004/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=V3DkORIO4Oe1mcf36ER0cQ%3D%3D&trackingId=rDmrWQuKXWGCa5naRgDVVA%3D%3D&trk=flagship3_search_srp_jobs"
$row6[0,11] = 236.33
$row6[0,12] = "plan, demand, planning, get, experience, channel, time, work, employee, market"
$row6[0,13] = "demand planning, forecasting, remote"
$row6[0,14] = "remote job, planning manager"
$row6[0,15] = "United States"
$row6[0,16] = "4 days ago"
$row6[0,17] = "2025-04-18"
$ws.Range("A6:R6").Value = $row6

# Row 7: hudsons - Director of Operations - Launch Stage
$row7 = New-Object "object[,]" 1,18
$row7[0,0] = "hudsons"
$row7[0,1] = "Director of Operations - Launch Stage"
$row7[0,2] = $false
$row7[0,3] = $false
$row7[0,4] = $false
$row7[0,5] = $true
$row7[0,6] = $false
$row7[0,7] = $true
$row7[0,8] = $false
$row7[0,9] = $false
$row7[0,10] = "https://www.linkedin.com/jobs/view/4285913612/?eBP=BUDGET_EXHAUSTED_JOB&refId=V3DkORIO4Oe1mcf36ER0cQ%3D%3D&trackingId=1z4Tro4WplzpX5b8m5%2FW%2Fg%3D%3D&trk=flagship3_search_srp_jobs"
$row7[0,11] = 239.02
$row7[0,12] = "apparel, operations, fashion, operational, build, systems, vendor, medical, mark, brand"
$row7[0,13] = "inventory planning, remote, remote work"
$row7[0,14] = "remote job, planning manager"
$row7[0,15] = "United States"
$row7[0,16] = "3 days ago"
$row7[0,17] = "2025-04-19"
$ws.Range("A7:R7").Value = $row7

# Row 8: UGG - Sr. Director, Demand Planning - UGG
$row8 = New-Object "object[,]" 1,18
$row8[0,0] = "UGG"
$row8[0,1] = "Sr. Director, Demand Planning - UGG"
$row8[0,2] = $false
$row8[0,3] = $false
$row8[0,4] = $false
$row8[0,5] = $true
$row8[0,6] = $false
$row8[0,7] = $true
$row8[0,8] = $false
$row8[0,9] = $false
$row8[0,10] = "https://www.linkedin.com/jobs/view/4281741341/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=5yq3jpTgHwZntidRX8LkfA%3D%3D&trackingId=YKnMvTyiBrMQM5LcKUDpJg%3D%3D&trk=flagship3_search_srp_jobs"
$row8[0,11] = 291
$row8[0,12] = "plan, demand, planning, get, experience, channel, time, work, employee, market"
$row8[0,13] = "demand planning, forecasting, remote"
$row8[0,14] = "remote job, planning manager"
$row8[0,15] = "United States"
$row8[0,16] = "4 days ago"
$row8[0,17] = "2025-04-18"
$ws.Range("A8:R8").Value = $row8

# Row 9: Intuvia Technologies - Supply Chain Analyst
$row9 = New-Object "object[,]" 1,18
$row9[0,0] = "Intuvia Technologies"
$row9[0,1] = "Supply Chain Analyst"
$row9[0,2] = $false
$row9[0,3] = $false
$row9[0,4] = $false
$row9[0,5] = $true
$row9[0,6] = $false
$row9[0,7] = $true
$row9[0,8] = $false
$row9[0,9] = $false
$row9[0,10] = "https://www.linkedin.com/jobs/view/4286802419/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=f40UJ2fSXPEN5Kq7MfUFDQ%3D%3D&trackingId=UYbadGpz7bXVTgE6AphuXQ%3D%3D&trk=flagship3_search_srp_jobs"
$row9[0,11] = 419.3
$row9[0,12] = "supply, manage, chain, management, product, demand, forecast, inventory, sales, skills"
$row9[0,13] = "demand planning, forecasting, remote"
$row9[0,14] = "remote job, planning manager"
$row9[0,15] = "United States"
$row9[0,16] = "1 day ago"
$row9[0,17] = "2025-04-21"
$ws.Range("A9:R9").Value = $row9

# Row 10: Insight Global - Capacity & Demand Analyst
$row10 = New-Object "object[,]" 1,18
$row10[0,0] = "Insight Global"
$row10[0,1] = "Capacity & Demand Analyst"
$row10[0,2] = $false
$row10[0,3] = $false
$row10[0,4] = $false
$row10[0,5] = $true
$row10[0,6] = $false
$row10[0,7] = $true
$row10[0,8] = $false
$row10[0,9] = $false
$row10[0,10] = "https://www.linkedin.com/jobs/view/4283516092/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=KPup3DMmYFpqExxaSUpLNg%3D%3D&trackingId=mzAtXW%2BPUoBZBS0YdM9mQw%3D%3D&trk=flagship3_search_srp_jobs"
$row10[0,11] = 467.32
$row10[0,12] = "experience, portfolio, capacity, plan, management, resource, skills, demand, analyst, project"
$row10[0,13] = "capacity planning, remote"
$row10[0,14] = "remote job, planning manager"
$row10[0,15] = "United States"
$row10[0,16] = "Capacity & Demand Analyst`nCapacity & Demand Analyst with verification`nInsight Global`nUnited States (Remote)`n`$40/hr - `$45/hr · Vision, 401(k)`nActively reviewing applicants`nViewed`nEasy Apply"
$row10[0,17] = $null
$ws.Range("A10:R10").Value = $row10

# Row 11: Spiraledge, Inc - Planning Analyst
$row11 = New-Object "object[,]" 1,18
$row11[0,0] = "Spiraledge, Inc"
$row11[0,1] = "Planning Analyst"
$row11[0,2] = $false
$row11[0,3] = $false
$row11[0,4] = $false
$row11[0,5] = $true
$row11[0,6] = $false
$row11[0,7] = $true
$row11[0,8] = $false
$row11[0,9] = $false
$row11[0,10] = "https://www.linkedin.com/jobs/view/4286502904/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=KPup3DMmYFpqExxaSUpLNg%3D%3D&trackingId=KM7KuvpHMYIu1uJ1Pep%2FPA%3D%3D&trk=flagship3_search_srp_jobs"
$row11[0,11] = 482.95
$row11[0,12] = "planning, per, work, com, team, demand, improve, inventory, sales, experience"
$row11[0,13] = "inventory planning, demand planning, forecasting, remote"
$row11[0,14] = "remote job, planning manager"
$row11[0,15] = "United States"
$row11[0,16] = "2 days ago"
$row11[0,17] = "2025-04-20"
$ws.Range("A11:R11").Value = $row11

# Row 12: Dynamic Events, Inc. - Event Manager - Logistics & Project Management Focus
$row12 = New-Object "object[,]" 1,18
$row12[0,0] = "Dynamic Events, Inc."
$row12[0,1] = "Event Manager - Logistics & Project Management Focus"
$row12[0,2] = $false
$row12[0,3] = $false
$row12[0,4] = $false
$row12[0,5] = $true
$row12[0,6] = $false
$row12[0,7] = $true
$row12[0,8] = $false
$row12[0,9] = $false
$row12[0,10] = "https://www.linkedin.com/jobs/view/4222843514/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=KPup3DMmYFpqExxaSUpLNg%3D%3D&trackingId=q%2BPDBugsrs%2FU8vwI1Mo4oA%3D%3D&trk=flagship3_search_srp_jobs"
$row12[0,11] = 491.23
$row12[0,12] = "manage, event, management, ability, project, team, over, time, events, client"
$row12[0,13] = "forecasting, remote, work from home, fully remote"
$row12[0,14] = "remote job, planning manager"
$row12[0,15] = "United States"
$row12[0,16] = "3 months ago"
$row12[0,17] = "2025-01-22"
$ws.Range("A12:R12").Value = $row12

# Row 13: UGG - Sr. Director, Demand Planning - UGG
$row13 = New-Object "object[,]" 1,18
$row13[0,0] = "UGG"
$row13[0,1] = "Sr. Director, Demand Planning - UGG"
$row13[0,2] = $false
$row13[0,3] = $false
$row13[0,4] = $false
$row13[0,5] = $true
$row13[0,6] = $false
$row13[0,7] = $true
$row13[0,8] = $false
$row13[0,9] = $false
$row13[0,10] = "https://www.linkedin.com/jobs/view/4281741342/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=KPup3DMmYFpqExxaSUpLNg%3D%3D&trackingId=uIDExrxtiFQgu0fXjRiSEg%3D%3D&trk=flagship3_search_srp_jobs"
$row13[0,11] = 544.4
$row13[0,12] = "plan, demand, planning, get, experience, channel, time, work, employee, market"
$row13[0,13] = "demand planning, forecasting, remote"
$row13[0,14] = "remote job, planning manager"
$row13[0,15] = "United States"
$row13[0,16] = "4 days ago"
$row13[0,17] = "2025-04-18"
$ws.Range("A13:R13").Value = $row13

# Row 14: Lensa - Sales & Operations Planning (S&OP) Supply Planner
$row14 = New-Object "object[,]" 1,18
$row14[0,0] = "Lensa"
$row14[0,1] = "Sales & Operations Planning (S&OP) Supply Planner"
$row14[0,2] = $true
$row14[0,3] = $false
$row14[0,4] = $false
$row14[0,5] = $true
$row14[0,6] = $false
$row14[0,7] = $true
$row14[0,8] = $false
$row14[0,9] = $false
$row14[0,10] = "https://www.linkedin.com/jobs/view/4287293367/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=KPup3DMmYFpqExxaSUpLNg%3D%3D&trackingId=WZMHxUhYbnfnKpULoyWHAw%3D%3D&trk=flagship3_search_srp_jobs"
$row14[0,11] = 558.17
$row14[0,12] = "supply, plan, com, job, age, lead, planning, vernova, process, business"
$row14[0,13] = "relocation assistance, supply planning, capacity planning, operations planning, forecasting, s&op, sales and operations planning, Relocation Assistance Provided: No, remote, remote position"
$row14[0,14] = "remote job, planning manager"
$row14[0,15] = "United States"
$row14[0,16] = "1 day ago"
$row14[0,17] = "2025-04-21"
$ws.Range("A14:R14").Value = $row14

# Row 15: Binding Minds Inc. (Certified Disability Owned Business Enterprise) - Finance and IT Project Manager
$row15 = New-Object "object[,]" 1,18
$row15[0,0] = "Binding Minds Inc. (Certified Disability Owned Business Enterprise)"
$row15[0,1] = "Finance and IT Project Manager"
$row15[0,2] = $false
$row15[0,3] = $false
$row15[0,4] = $false
$row15[0,5] = $true
$row15[0,6] = $false
$row15[0,7] = $true
$row15[0,8] = $false
$row15[0,9] = $false
$row15[0,10] = "https://www.linkedin.com/jobs/view/4284517016/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=KPup3DMmYFpqExxaSUpLNg%3D%3D&trackingId=TEJ8xRrj9YChv%2BEPzskw5w%3D%3D&trk=flagship3_search_srp_jobs"
$row15[0,11] = 566.04
$row15[0,12] = "project, manage, management, finance, implement, experience, financial, solution, implementation, team"
$row15[0,13] = "forecasting, remote"
$row15[0,14] = "remote job, planning manager"
$row15[0,15] = "United States"
$row15[0,16] = "3 days ago"
$row15[0,17] = "2025-04-19"
$ws.Range("A15:R15").Value = $row15

# Row 16: Kate McLeod - Senior Manager, Planning & Operations
$row16 = New-Object "object[,]" 1,18
$row16[0,0] = "Kate McLeod"
$row16[0,1] = "Senior Manager, Planning & Operations"
$row16[0,2] = $false
$row16[0,3] = $false
$row16[0,4] = $false
$row16[0,5] = $true
$row16[0,6] = $false
$row16[0,7] = $true
$row16[0,8] = $false
$row16[0,9] = $false
$row16[0,10] = "https://www.linkedin.com/jobs/view/4197615128/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=hZpYSk3tcglwNmk6OuehTQ%3D%3D&trackingId=6D8niZEI%2B%2Bm3iACClrVnVg%3D%3D&trk=flagship3_search_srp_jobs"
$row16[0,11] = 717.75
$row16[0,12] = "per, plan, age, manage, planning, product, operations, demand, production, supply"
$row16[0,13] = "production planning, demand planning, supply chain planning, forecasting, remote"
$row16[0,14] = "remote job, planning manager"
$row16[0,15] = "United States"
$row16[0,16] = "4 months ago"
$row16[0,17] = "2024-12-23"
$ws.Range("A16:R16").Value = $row16

# Row 17: Anker Innovations LTD - Supply Chain Manager
$row17 = New-Object "object[,]" 1,18
$row17[0,0] = "Anker Innovations LTD"
$row17[0,1] = "Supply Chain Manager"
$row17[0,2] = $false
$row17[0,3] = $false
$row17[0,4] = $false
$row17[0,5] = $true
$row17[0,6] = $false
$row17[0,7] = $true
$row17[0,8] = $false
$row17[0,9] = $false
$row17[0,10] = "https://www.linkedin.com/jobs/view/4278650679/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=hZpYSk3tcglwNmk6OuehTQ%3D%3D&trackingId=LiQIfQUJQ2wnoNlz6DuX6w%3D%3D&trk=flagship3_search_srp_jobs"
$row17[0,11] = 746.47
$row17[0,12] = "com, age, anker, plan, product, planning, demand, work, chain, forecast"
$row17[0,13] = "production planning, demand planning, supply chain planning, forecasting, remote"
$row17[0,14] = "remote job, planning manager"
$row17[0,15] = "United States"
$row17[0,16] = "1 day ago"
$row17[0,17] = "2025-04-21"
$ws.Range("A17:R17").Value = $row17

# Row 18: GlobalSource IT - Supply Chain Transformation Project Manager
$row18 = New-Object "object[,]" 1,18
$row18[0,0] = "GlobalSource IT"
$row18[0,1] = "Supply Chain Transformation Project Manager"
$row18[0,2] = $false
$row18[0,3] = $false
$row18[0,4] = $false
$row18[0,5] = $true
$row18[0,6] = $false
$row18[0,7] = $true
$row18[0,8] = $false
$row18[0,9] = $false
$row18[0,10] = "https://www.linkedin.com/jobs/view/4278321337/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=DS3KwMrj3cvtJXyI0Ch06w%3D%3D&trackingId=Xw%2BL%2FxzgjP%2B85BxwArXg%2BQ%3D%3D&trk=flagship3_search_srp_jobs"
$row18[0,11] = 803.03
$row18[0,12] = "supply, chain, project, function, able, functional, erp, operations, lead, transformation"
$row18[0,13] = "supply planning, demand planning, operations planning, s&op, sales and operations planning, remote"
$row18[0,14] = "remote job, planning manager"
$row18[0,15] = "United States"
$row18[0,16] = "1 week ago"
$row18[0,17] = "2025-04-15"
$ws.Range("A18:R18").Value = $row18

# Row 19: Lantern - Strategy Analyst
$row19 = New-Object "object[,]" 1,18
$row19[0,0] = "Lantern"
$row19[0,1] = "Strategy Analyst"
$row19[0,2] = $false
$row19[0,3] = $false
$row19[0,4] = $false
$row19[0,5] = $true
$row19[0,6] = $false
$row19[0,7] = $true
$row19[0,8] = $false
$row19[0,9] = $false
$row19[0,10] = "https://www.linkedin.com/jobs/view/4266378939/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=xcCHYviEbmOQ6NTfxgeSJQ%3D%3D&trackingId=%2FH%2FTdfesqeCEtbCcd1O1HQ%3D%3D&trk=flagship3_search_srp_jobs"
$row19[0,11] = 987.42
$row19[0,12] = "care, lantern, team, com, support, high, out, product, drive, opportunities"
$row19[0,13] = "forecasting, remote"
$row19[0,14] = "remote job, planning manager"
$row19[0,15] = "United States"
$row19[0,16] = "Strategy Analyst`nStrategy Analyst with verification`nLantern`nUnited States (Remote)`nMedical, Dental, Vision, 401(k)`nActively reviewing applicants`nViewed`nEasy Apply"
$row19[0,17] = $null
$ws.Range("A19:R19").Value = $row19

# Row 20: Lyra Health - Finance Manager - FP&A
$row20 = New-Object "object[,]" 1,18
$row20[0,0] = "Lyra Health"
$row20[0,1] = "Finance Manager - FP&A"
$row20[0,2] = $false
$row20[0,3] = $false
$row20[0,4] = $false
$row20[0,5] = $true
$row20[0,6] = $false
$row20[0,7] = $true
$row20[0,8] = $false
$row20[0,9] = $false
$row20[0,10] = "https://www.linkedin.com/jobs/view/4286556685/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=xcCHYviEbmOQ6NTfxgeSJQ%3D%3D&trackingId=UgL9qYbYHK5uCETYHdeKOw%3D%3D&trk=flagship3_search_srp_jobs"
$row20[0,11] = 1008.36
$row20[0,12] = "form, per, format, information, com, financial, lyra, out, inc, health"
$row20[0,13] = "forecasting, remote, fully remote"
$row20[0,14] = "remote job, planning manager"
$row20[0,15] = "United States"
$row20[0,16] = "2 days ago"
$row20[0,17] = "2025-04-20"
$ws.Range("A20:R20").Value = $row20

# Row 21: Insight Global - Business Analyst
$row21 = New-Object "object[,]" 1,18
$row21[0,0] = "Insight Global"
$row21[0,1] = "Business Analyst"
$row21[0,2] = $false
$row21[0,3] = $false
$row21[0,4] = $false
$row21[0,5] = $true
$row21[0,6] = $false
$row21[0,7] = $true
$row21[0,8] = $false
$row21[0,9] = $false
$row21[0,10] = "https://www.linkedin.com/jobs/view/4286156458/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=gfolC6em6DL3JYEwja3GRw%3D%3D&trackingId=MgpQkZ%2Bi6bgVVMrXBDrO3A%3D%3D&trk=flagship3_search_srp_jobs"
$row21[0,11] = 1409.53
$row21[0,12] = "data, business, support, operations, performance, experience, work, operational, ability, team"
$row21[0,13] = "forecasting, remote, fully remote"
$row21[0,14] = "remote job, planning manager"
$row21[0,15] = "United States"
$row21[0,16] = "2 days ago"
$row21[0,17] = "2025-04-20"
$ws.Range("A21:R21").Value = $row21

# Row 22: Aquent - Senior Manager Strategic Planning
$row22 = New-Object "object[,]" 1,18
$row22[0,0] = "Aquent"
$row22[0,1] = "Senior Manager Strategic Planning"
$row22[0,2] = $false
$row22[0,3] = $false
$row22[0,4] = $false
$row22[0,5] = $true
$row22[0,6] = $false
$row22[0,7] = $true
$row22[0,8] = $false
$row22[0,9] = $false
$row22[0,10] = "https://www.linkedin.com/jobs/view/4281390111/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=gfolC6em6DL3JYEwja3GRw%3D%3D&trackingId=1Y%2BBgeLomA8Sp2Ggf453Ew%3D%3D&trk=flagship3_search_srp_jobs"
$row22[0,11] = 1443.39
$row22[0,12] = "plan, plannin, rate, planning, levels, product, senior, forecast, manager, low"
$row22[0,13] = "demand planning, forecasting, remote"
$row22[0,14] = "remote job, planning manager"
$row22[0,15] = "United States"
$row22[0,16] = "1 week ago"
$row22[0,17] = "2025-04-15"
$ws.Range("A22:R22").Value = $row22

# Row 23: Alliant Human Capital - Production and Supply Chain Manager
$row23 = New-Object "object[,]" 1,18
$row23[0,0] = "Alliant Human Capital"
$row23[0,1] = "Production and Supply Chain Manager"
$row23[0,2] = $false
$row23[0,3] = $false
$row23[0,4] = $false
$row23[0,5] = $true
$row23[0,6] = $false
$row23[0,7] = $true
$row23[0,8] = $false
$row23[0,9] = $false
$row23[0,10] = "https://www.linkedin.com/jobs/view/4236409944/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=gtDkp%2FpCKTdB5YdrxqZDbw%3D%3D&trackingId=p4vvAt0t1vop6%2FgXzxzYyg%3D%3D&trk=flagship3_search_srp_jobs"
$row23[0,11] = 1493.33
$row23[0,12] = "plan, manage, product, planning, production, demand, operations, supply, chain, management"
$row23[0,13] = "production planning, demand planning, supply chain planning, forecasting, remote"
$row23[0,14] = "remote job, planning manager"
$row23[0,15] = "United States"
$row23[0,16] = "3 weeks ago"
$row23[0,17] = "2025-04-01"
$ws.Range("A23:R23").Value = $row23

# Row 24: NSI Group, LLC - Operations Planning Manager
$row24 = New-Object "object[,]" 1,18
$row24[0,0] = "NSI Group, LLC"
$row24[0,1] = "Operations Planning Manager"
$row24[0,2] = $false
$row24[0,3] = $false
$row24[0,4] = $false
$row24[0,5] = $true
$row24[0,6] = $false
$row24[0,7] = $true
$row24[0,8] = $false
$row24[0,9] = $false
$row24[0,10] = "https://www.linkedin.com/jobs/view/4121212510/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=gtDkp%2FpCKTdB5YdrxqZDbw%3D%3D&trackingId=Tj7zALoquCd%2F9tuojCNgWA%3D%3D&trk=flagship3_search_srp_jobs"
$row24[0,11] = 1512.89
$row24[0,12] = "manage, plan, product, planning, quality, operations, food, nsi, management, production"
$row24[0,13] = "production planning, operations planning, forecasting, remote, remote work, remote position"
$row24[0,14] = "remote job, planning manager"
$row24[0,15] = "United States"
$row24[0,16] = "2 weeks ago"
$row24[0,17] = "2025-04-08"
$ws.Range("A24:R24").Value = $row24

# --- Drop the two trailing rows that no longer exist in the refreshed data ---
$ws.Range("A25:A26").EntireRow.Delete()

Write-Host "Edit complete"